$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for years 2010-2022 (replaces old 2000-2020 data)
$data = @(
    @("2010年", 89874.16, 17.8, 83101.50999999999, 21.3),
    @("2011年", 109247.79, 21.6, 103874.43, 25),
    @("2012年", 125952.97, 15.3, 117253.52, 12.9),
    @("2013年", 140212.1, 11.3, 129209.64, 10.2),
    @("2014年", 151785.56, 8.300000000000001, 140370.03, 8.6),
    @("2015年", 175877.77, 13.2, 152269.23, 5.8),
    @("2016年", 187755.21, 6.3, 159604.97, 4.5),
    @("2017年", 203085.49, 7.6, 172592.77, 7.4),
    @("2018年", 220904.13, 8.699999999999999, 183359.84, 6.2),
    @("2019年", 238858.37, 8.1, 190390.08, 3.8),
    @("2020年", 245679.03, 2.9, 182913.88, -3.9),
    @("2021年", 245673, 0, 202554.64, 10.7),
    @("2022年", 260609.17, 6.1, 203703.48, 0.6)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Remove the now-unused rows that previously held data through row 22
$lastOldRow = 22
$newLastRow = $startRow + $data.Count - 1
if ($lastOldRow -gt $newLastRow) {
    $deleteRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($lastOldRow, 5))
    $deleteRange.EntireRow.Delete()
}
